$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "True:" / "False:" / "False," lead-in from the explanation cells
# (columns I, J, K) across the four question rows.
$ws.Range("I2").Value = 'The term inflation shock is used to refer to an exogenous shift in the Phillips curve.'
$ws.Range("J2").Value = 'The central bank needs to increase the interest rate immediately to decrease output below equilibrium.'
$ws.Range("K2").Value = 'Inflation decreases if the nominal interest rate is above the stabilizing interest rate. It could be that the interest rate is increasing, but it is still below the stabilizing interest rate. Then inflation increases.'

$ws.Range("I3").Value = 'If the central bank does not change the interest rate, inflation soars up. Hence, it chooses a not sustainable policy.'
$ws.Range("J3").Value = 'A more inflation-averse central bank’s immediate response is to choose a higher interest rate than a less inflation-averse central bank. A higher increase in the interest rate leads to a larger decrease in the output.'
$ws.Range("K3").Value = 'In order to shorten the period of inflation away from target a more inflation-averse central bank increases the interest rates greater than a less inflation-averse one.'

$ws.Range("I4").Value = 'The value of BETA does not reflect whether the central bank focuses on achieving an inflation target or an output target. Rather, a central bank with a lower BETA is willing to trade off a longer period during which inflation is away from target to reduce the impact on unemployment of the adjustment path back to equilibrium than would a more inflation-averse central bank with a higher BETA.'
$ws.Range("J4").Value = 'BETA < 1 characterizes a central bank that places a greater weight on deviations in inflation target than on deviations in employment.'
$ws.Range("K4").Value = 'The central bank is indifferent between inflation and output deviations from targets.'

$ws.Range("I5").Value = 'If the expectations are better formed and the central bank has a high level of credibility, it needs to increase the interest rate less. This leads to a lower decrease in output and a lower increase in unemployment than a less credible central bank.'
$ws.Range("J5").Value = 'BETA does not reflect the credibility of the central bank.'
$ws.Range("K5").Value = 'If the public is better informed about monetary policy, and the central bank has a high level of credibility, it is more likely that inflation expectations will stay close to the target.'

# D3/E3 carry redundant formatting (a font style identical to the default) -
# normalise them back onto the plain wrap/top-aligned style used elsewhere.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").VerticalAlignment = -4160
$ws.Range("D3").WrapText = $true
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").VerticalAlignment = -4160
$ws.Range("E3").WrapText = $true

# Row 4 shrinks now that its explanations are shorter.
$ws.Rows.Item(4).RowHeight = 174

# Selection moved on to J5 for the next edit session.
[void]$ws.Range("J5").Select()
